$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.275.46'
$ws.Range('E2').Value = '  -3.72%  '
$ws.Range('D3').Value = '3.150.61'
$ws.Range('E3').Value = '  -3.34%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '607.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.33'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.89%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '3.148.52'
$ws.Range('E8').Value = '  -3.37%  '
$ws.Range('E9').Value = '  -4.69%  '
$ws.Range('E10').Value = '  -6.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.50'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -6.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.476'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -6.00%  '
$ws.Range('E13').Value = '  -7.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.62'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -9.28%  '
$ws.Range('D15').Value = '3.671.93'
$ws.Range('E15').Value = '  -3.27%  '
$ws.Range('D16').Value = '64.288.19'
$ws.Range('E16').Value = '  -3.78%  '
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('D18').Value = '3.152.86'
$ws.Range('E18').Value = '  -3.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.92'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '480.51'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.70%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.77'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.48%  '
$ws.Range('E22').Value = '  -5.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.78'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.67'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -8.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.57'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('E27').Value = '  -5.19%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.46'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.35%  '
$ws.Range('E29').Value = '  -9.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.81'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.00%  '
$ws.Range('E31').Value = '  -18.82%  '
$ws.Range('E32').Value = '  -5.98%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.21'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.82%  '
$ws.Range('E35').Value = '  -4.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '54.33'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.96'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -8.12%  '
$ws.Range('D38').Value = '0.0₃0730'
$ws.Range('E38').Value = '  -8.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '455.61'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.02%  '
$ws.Range('E40').Value = '  -14.20%  '
$ws.Range('E41').Value = '  -7.38%  '
$ws.Range('E42').Value = '  -8.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.43'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.91%  '
$ws.Range('D44').Value = '2.848.57'
$ws.Range('E44').Value = '  -4.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.265'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -9.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.26'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -10.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.60'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.99%  '
$ws.Range('E49').Value = '  -7.07%  '
$ws.Range('E50').Value = '  -4.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '119.69'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.55%  '
